{"js": "const replacements = [\n  [\"2024-03-02 Saturday\", \"2024-03-03 Sunday\"],\n  [\"112\u00f78=14, 0\", \"561\u00f72=280, 1\"],\n  [\"428\u00f79=47, 5\", \"828\u00f72=414, 0\"],\n  [\"289\u00f78=36, 1\", \"292\u00f75=58, 2\"],\n  [\"262\u00f74=65, 2\", \"751\u00f77=107, 2\"],\n  [\"110\u00f76=18, 2\", \"799\u00f72=399, 1\"],\n  [\"602\u00f77=86, 0\", \"572\u00f74=143, 0\"],\n  [\"706\u00f73=235, 1\", \"733\u00f79=81, 4\"],\n  [\"106\u00f72=53, 0\", \"540\u00f72=270, 0\"],\n  [\"785\u00f74=196, 1\", \"612\u00f79=68, 0\"],\n  [\"995\u00f72=497, 1\", \"566\u00f76=94, 2\"],\n  [\"113\u00f78=14, 1\", \"151\u00f72=75, 1\"],\n  [\"234\u00f78=29, 2\", \"356\u00f73=118, 2\"],\n  [\"429\u00f75=85, 4\", \"951\u00f78=118, 7\"],\n  [\"224\u00f75=44, 4\", \"241\u00f75=48, 1\"],\n  [\"266\u00f77=38, 0\", \"212\u00f77=30, 2\"],\n  [\"851\u00f79=94, 5\", \"501\u00f76=83, 3\"],\n  [\"556\u00f77=79, 3\", \"474\u00f73=158, 0\"],\n  [\"154\u00f74=38, 2\", \"586\u00f73=195, 1\"],\n  [\"181\u00f72=90, 1\", \"600\u00f79=66, 6\"],\n  [\"811\u00f74=202, 3\", \"305\u00f76=50, 5\"],\n  [\"197\u00f76=32, 5\", \"346\u00f73=115, 1\"],\n  [\"226\u00f72=113, 0\", \"933\u00f74=233, 1\"],\n  [\"421\u00f75=84, 1\", \"180\u00f72=90, 0\"],\n  [\"378\u00f75=75, 3\", \"180\u00f76=30, 0\"],\n  [\"118\u00f77=16, 6\", \"662\u00f79=73, 5\"],\n];\n\nconst body = context.document.body;\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Not found: ' + from);\n  }\n  for (const item of results.items) {\n    item.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nreturn 'ok';", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-03-02 Saturday\", \"2024-03-03 Sunday\"),\n  @(\"112\u00f78=14, 0\", \"561\u00f72=280, 1\"),\n  @(\"428\u00f79=47, 5\", \"828\u00f72=414, 0\"),\n  @(\"289\u00f78=36, 1\", \"292\u00f75=58, 2\"),\n  @(\"262\u00f74=65, 2\", \"751\u00f77=107, 2\"),\n  @(\"110\u00f76=18, 2\", \"799\u00f72=399, 1\"),\n  @(\"602\u00f77=86, 0\", \"572\u00f74=143, 0\"),\n  @(\"706\u00f73=235, 1\", \"733\u00f79=81, 4\"),\n  @(\"106\u00f72=53, 0\", \"540\u00f72=270, 0\"),\n  @(\"785\u00f74=196, 1\", \"612\u00f79=68, 0\"),\n  @(\"995\u00f72=497, 1\", \"566\u00f76=94, 2\"),\n  @(\"113\u00f78=14, 1\", \"151\u00f72=75, 1\"),\n  @(\"234\u00f78=29, 2\", \"356\u00f73=118, 2\"),\n  @(\"429\u00f75=85, 4\", \"951\u00f78=118, 7\"),\n  @(\"224\u00f75=44, 4\", \"241\u00f75=48, 1\"),\n  @(\"266\u00f77=38, 0\", \"212\u00f77=30, 2\"),\n  @(\"851\u00f79=94, 5\", \"501\u00f76=83, 3\"),\n  @(\"556\u00f77=79, 3\", \"474\u00f73=158, 0\"),\n  @(\"154\u00f74=38, 2\", \"586\u00f73=195, 1\"),\n  @(\"181\u00f72=90, 1\", \"600\u00f79=66, 6\"),\n  @(\"811\u00f74=202, 3\", \"305\u00f76=50, 5\"),\n  @(\"197\u00f76=32, 5\", \"346\u00f73=115, 1\"),\n  @(\"226\u00f72=113, 0\", \"933\u00f74=233, 1\"),\n  @(\"421\u00f75=84, 1\", \"180\u00f72=90, 0\"),\n  @(\"378\u00f75=75, 3\", \"180\u00f76=30, 0\"),\n  @(\"118\u00f77=16, 6\", \"662\u00f79=73, 5\"),\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $found) {\n    throw \"Not found: \" + $pair[0]\n  }\n}\n\nWrite-Output \"done\""}
